$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 2, shifting existing rows (2..179) down to (3..180)
$ws.Rows("2:2").Insert()

# Populate the newly inserted row 2 with values (mirrors old row 2's
# Transaction Type / Payment Type / InternalComment, new USD Amount)
$ws.Range("E2").Value = "Deposit"
$ws.Range("N2").Value = "Crypto"
$ws.Range("P2").Value = "ETH"
$ws.Range("T2").Value = 2123.5902000000001
$ws.Range("AB2").Style = $ws.Range("AB3").Style

# Update sheet view: topLeftCell / selection
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("E2:P3").Select()
$excel.ActiveCell = $ws.Range("P3")
